$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 9 (shifts old rows 9-13 down to 10-14)
$ws.Rows.Item(9).Insert()

# Make sure the newly-inserted blank row actually materializes its (empty) cells
$ws.Range("A9:D9").Style = "Normal"

# Fill in row 8 (previously blank spacer row) with a new "DK/NA" entry
$ws.Cells.Item(8, 2).Value = "DK/NA"
$ws.Cells.Item(8, 3).Value = "(0.0%) 0"
$ws.Cells.Item(8, 4).Value = "(0.0%) 0"

# Add a new "DK/NA" row at the end of the Gender block (row 15)
$ws.Range("A15").Style = "Normal"
$ws.Cells.Item(15, 2).Value = "DK/NA"
$ws.Cells.Item(15, 3).Value = "(0.0%) 0"
$ws.Cells.Item(15, 4).Value = "(0.0%) 0"
